# Auto update Excel log
# Appends new sensor-log rows (2026-01-28, ~15:15-15:16) exported from the
# SeniorConnect sensors to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

function Add-LogRows($Worksheet, $StartRow, $Rows) {
    $endRow = $StartRow + $Rows.Length - 1
    $rng = $Worksheet.Range("A" + $StartRow + ":F" + $endRow)
    # Force text formatting first so values like "2026-01-28" and "15:00"
    # are stored as literal strings instead of being auto-converted to
    # date/time values by Excel's smart input parsing.
    $rng.NumberFormat = "@"

    for ($i = 0; $i -lt $Rows.Length; $i++) {
        $r = $StartRow + $i
        $rowVals = $Rows[$i]
        for ($j = 0; $j -lt $rowVals.Length; $j++) {
            $col = $j + 1
            $Worksheet.Cells.Item($r, $col).Value = $rowVals[$j]
        }
    }

    # Restore default (General) formatting now that the literal text is
    # already committed to the cells, so no stray number-format style is
    # left behind on the new rows.
    $rng.ClearFormats()
}

# --- PIR sheet: rows 281-293 ---
$pirRows = @(
    @("2026-01-28","15:15:32","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:15:33","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:15:38","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:15:43","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:15:49","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:15:53","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:15:58","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:16:03","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:16:09","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:16:13","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:16:18","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:16:23","15:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","15:16:28","15:00","Bathroom","No Motion","Inactive")
)
$wsPIR = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPIR 281 $pirRows

# --- Humidity sheet: rows 268-279 ---
$humidityRows = @(
    @("2026-01-28","15:15:30","15:00","Bathroom","87.4%","Active"),
    @("2026-01-28","15:15:40","15:00","Bathroom","87.4%","Active"),
    @("2026-01-28","15:15:44","15:00","Bathroom","88.4%","Active"),
    @("2026-01-28","15:15:48","15:00","Bathroom","87.4%","Active"),
    @("2026-01-28","15:15:52","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:16:00","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:16:04","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:16:08","15:00","Bathroom","87.4%","Active"),
    @("2026-01-28","15:16:12","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:16:16","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:16:20","15:00","Bathroom","88.3%","Active"),
    @("2026-01-28","15:16:24","15:00","Bathroom","88.4%","Active")
)
$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHumidity 268 $humidityRows

# --- Temperature sheet: rows 268-279 ---
$temperatureRows = @(
    @("2026-01-28","15:15:31","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:15:40","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:15:45","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:15:48","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:15:52","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:16:01","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:16:05","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:16:09","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:16:13","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:16:17","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:16:21","15:00","Bathroom","22.9C","Active"),
    @("2026-01-28","15:16:25","15:00","Bathroom","22.9C","Active")
)
$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRows $wsTemperature 268 $temperatureRows
